$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (before old row 4),
# shifting all existing data rows down by two positions.
$ws.Rows("4:5").Insert()

# Populate the two newly inserted rows with the latest weekly price data.
# Row 4: Uva - Red Globe
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44631
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100109
$ws.Range("H4").Value = "Uva"
$ws.Range("I4").Value = 100109001
$ws.Range("J4").Value = "Uva"
$ws.Range("K4").Value = "Red Globe"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("Q4").Value = "`$/caja 20 kilos"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 1225
$ws.Range("T4").Value = 20

# Row 5: Uva - Thompson seedless
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44631
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100109
$ws.Range("H5").Value = "Uva"
$ws.Range("I5").Value = 100109001
$ws.Range("J5").Value = "Uva"
$ws.Range("K5").Value = "Thompson seedless"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("Q5").Value = "`$/caja 20 kilos"
$ws.Range("R5").Value = "Región de Coquimbo"
$ws.Range("S5").Value = 1225
$ws.Range("T5").Value = 20
